$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'309.08"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-0.58%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'37.15"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'-2.59%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.129"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'0.30%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.07836"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'-1.00%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'4.406"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'0.20%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'8.274"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'0.44%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'1.875"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-1.25%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D10").Value = "'0.9255"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'-0.18%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.1175"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'-2.03%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.1894"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-0.50%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.08870"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-4.28%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.03315"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'-1.85%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.09603"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'0.10%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.001379"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'1.16%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.006205"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'6.60%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'3.390"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'-3.88%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.3457"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'0.29%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'6.394"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'21.12%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.1292"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'0.89%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D23").Value = "'0.04349"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'-0.53%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.001200"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'-3.86%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.004272"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'-0.06%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.0001402"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'7.78%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'0.0002904"
$ws.Range("D27").Style = "Normal"
$ws.Range("D39").Value = "'0.02156"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'2.54%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.04995"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'-1.56%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.007563"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'-1.27%"
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'0.30%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.008521"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-6.44%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.002013"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-0.40%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.008792"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'1.82%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006572"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'-1.51%"
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'0.14%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.003297"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'13.65%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.001445"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'20.79%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.00002103"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'0.14%"
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.0002003"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'0.14%"
$ws.Range("E51").Style = "Normal"
